$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.404.92'
$ws.Range("E2").Value = '  -0.18%  '

$ws.Range("D3").Value = '2.645.33'
$ws.Range("E3").Value = '  -0.14%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.06%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.76%  '

$ws.Range("E9").Value = '  +0.44%  '

$ws.Range("E10").Value = '  -0.95%  '

$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.366'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.80%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.152'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.32'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.62%  '

$ws.Range("D14").Value = '3.120.58'
$ws.Range("E14").Value = '  -0.32%  '

$ws.Range("D15").Value = '63.275.74'
$ws.Range("E15").Value = '  -0.30%  '

$ws.Range("E16").Value = '  -1.12%  '

$ws.Range("D17").Value = '2.640.72'
$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.45'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '342.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.18%  '

$ws.Range("E21").Value = '  +2.47%  '

$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.34%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.57'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.59%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.13%  '

$ws.Range("B27").Value = 'SuiNetwork'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.49%  '

$ws.Range("B28").Value = 'Bittensor'
$ws.Range("C28").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '548.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.43%  '

$ws.Range("E29").Value = '  -2.35%  '

$ws.Range("E30").Value = '  -0.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.06'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.41%  '

$ws.Range("E33").Value = '  -2.49%  '

$ws.Range("E34").Value = '  -1.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '167.16'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.58%  '

$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.406'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.70%  '

$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.05'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.89'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.11%  '

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '169.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.41%  '

$ws.Range("E44").Value = '  +1.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0578'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.74%  '

$ws.Range("E46").Value = '  -0.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0245'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.56%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0962'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.63%  '

$ws.Range("E50").Value = '  +3.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.73%  '
